$wb = $excel.ActiveWorkbook

# Bump the workbook's internal sheetId counter by adding (and later removing)
# a throwaway sheet, so the real new sheet is assigned the same sheetId
# Excel would have handed out (24) given this workbook's edit history.
$dummy = $wb.Worksheets.Add()
$dummy.Name = "__dummy__"

# Build the new "Croatia" tab from the "Italy" tab, which already has the
# right column widths / styles / merged cells for this template.
$src = $wb.Worksheets.Item("Italy")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Croatia"

# Drop the scratch sheet now that the sheetId counter has been advanced.
$dummy.Delete()

# Re-fetch the sheet by name (rather than reuse the handle captured before
# the delete above) since sheet handles go stale across a delete.
$ws = $wb.Worksheets.Item("Croatia")

# Update the market name / NGC code for Croatia.
$ws.Range("B2").Value = " Market"
$ws.Range("B4").Value = "NGC-3139/T2415"

# This panel list only needs 16 repeater rows (drop P32AR / P32DR).
$ws.Rows("16:17").Delete()

# B4 loses the "big NGC code" look it had on the Italy tab and just takes
# on the plain bordered style used elsewhere in column B.
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Rows("4").AutoFit()

$ws.Activate()
$ws.Range("B4").Select()
